$d = $word.ActiveDocument

# Replace the placeholder ID text (this also absorbs the trailing
# space run, so only a single run remains afterwards).
$d.Content.Find.Execute("**ID__AFFARS_mp_5306_502_topic_2__ID** ", $true, $false, $false, $false, $false, $true, 1, $false, "**ID__AFFARS_MP_5306_502_1__ID**", 2)

# Update the first paragraph's formatting: add a (invisible/no-line)
# paragraph border with 5pt spacing on all sides, and widen the left
# indent from 120 to 225 twips.
$para = $d.Paragraphs(1)
$para.Range.Borders.DistanceFromTop = 5
$para.Range.Borders.DistanceFromLeft = 5
$para.Range.Borders.DistanceFromBottom = 5
$para.Range.Borders.DistanceFromRight = 5
$para.Format.LeftIndent = 11.25
